# "admin user ngga bisa import" - rework the import template header row:
# drop status/id_level/timestamp/token columns, add id_role, rename layout,
# and clear the old legend/notes block in columns K:L.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): E1 "status" -> "id_role", F1 "id_level" -> "username";
# old G1/H1/I1 (timestamp/username/token) become blank placeholder header cells.
$ws.Range("E1").Value = "id_role"
$ws.Range("F1").Value = "username"
$ws.Range("G1").ClearContents()
$ws.Range("H1").ClearContents()
$ws.Range("I1").ClearContents()

# Row 2 sample/helper values in F2:I2 are no longer needed.
$ws.Range("F2:I2").ClearContents()

# Remove the old "keterangan" legend block (K2:L3) and the stray K5 "level" cell.
$ws.Range("K2:L3").ClearContents()
$ws.Range("K5").ClearContents()

# Move the active selection to where the old K5 cell used to be (now empty).
$ws.Range("I3").Select()
